$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.435222387313843
$ws.Range("B1").Value = 1.995905518531799
$ws.Range("C1").Value = 2.456494331359863
$ws.Range("D1").Value = 4.804220199584961
$ws.Range("E1").Value = 0.8637256622314453
